# Fixed a bug in Respin
#
# The rows in Sheet1 (A2:F23) were re-ordered (a permutation of whole
# rows - each row keeps its A..F values together but moves to a
# different row position). Rows 24, 25 and 26 are untouched.
#
# Read every source row into memory first, then write each row's
# values into its new destination row, so that overlapping
# source/destination ranges don't clobber data that still needs to be
# read.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) values for rows 2-23, columns A-F (1-6).
$data = @{}
for ($row = 2; $row -le 23; $row++) {
    $rowvals = @()
    for ($col = 1; $col -le 6; $col++) {
        $rowvals += $ws.Cells.Item($row, $col).Value2
    }
    $data[$row] = $rowvals
}

# Mapping of source row (key) -> destination row (value).
$map = @{
    2 = 15
    3 = 12
    4 = 8
    5 = 3
    6 = 14
    7 = 5
    8 = 6
    9 = 7
    10 = 13
    11 = 9
    12 = 10
    13 = 4
    14 = 11
    15 = 2
    16 = 18
    17 = 21
    18 = 17
    19 = 16
    20 = 19
    21 = 20
    22 = 23
    23 = 22
}

foreach ($srcRow in $map.Keys) {
    $destRow = $map[$srcRow]
    $vals = $data[$srcRow]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($destRow, $col).Value = $vals[$col - 1]
    }
}
